$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Kupon ID value in O2 (KUPON_ID column) from KPN20220000009 -> KPN20220000008
$ws.Range("O2").Value = "KPN20220000008"

# Update the PREPARATION text in F2 to reference the new Kupon ID as well
$ws.Range("F2").Value = "Username : 37841;`nPassword : bni1234;`nRole : RL09 Penyelia Settlement;`nKupon ID : KPN20220000008"

# Update the active selection in the sheet view to G2 (matches the diff's sheetView selection change)
[void]$ws.Range("G2").Select()
